$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.4205231666564941
$ws.Range("B1").Value = 2.004952430725098
$ws.Range("C1").Value = 2.143982887268066
$ws.Range("D1").Value = 1.924914717674255
$ws.Range("E1").Value = 0.9912518262863159
